$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.485.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.66%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.805.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.85%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'1.008"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'308.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.45%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.07117"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8753"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.89%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07743"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.96%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'19.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.841.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.80%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.97%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.39%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'85.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.89%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.86%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008569"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'26.530.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -3.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.970"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.94%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.87%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.976"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.18%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'150.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.91%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'17.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.27%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'112.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.824"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08653"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.91%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.044"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.58%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.56%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.428"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.008"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.79%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.564"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.22%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01928"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05087"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.878"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.934"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.90%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1564"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.103"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.29%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.009"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.83%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.96%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'9.936"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.38%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Quant"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'101.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.98%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.584"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.56%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'63.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.99%  "
$ws.Range("E51").Style = "Normal"
